$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1463.5385
$ws.Range("I96").Value = 1006
$ws.Range("J96").Value = 2195.6
$ws.Range("K96").Value = 3018
$ws.Range("L96").Value = 6586.799999999999
$ws.Range("M96").Value = -1645
$ws.Range("N96").Value = -9332.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 265366.72
$ws.Range("I61").Value = 205573
$ws.Range("K61").Value = 205573
$ws.Range("M61").Value = -205361
$ws.Range("H136").Value = 265366.72
$ws.Range("I136").Value = 205573
$ws.Range("K136").Value = 616719
$ws.Range("M136").Value = -614169

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3821.5186
$ws.Range("J94").Value = 7176.385
$ws.Range("L94").Value = 7176.385
$ws.Range("N94").Value = -8078.385
$ws.Range("H132").Value = 2525.2964
$ws.Range("I132").Value = 1101
$ws.Range("K132").Value = 3303
$ws.Range("M132").Value = -773

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1576.1111
$ws.Range("J5").Value = 2173.3333
$ws.Range("L5").Value = 6519.999899999999
$ws.Range("N5").Value = -6743.999899999999
$ws.Range("H80").Value = 866.6667
$ws.Range("J80").Value = 866.6667
$ws.Range("L80").Value = 2600.0001
$ws.Range("N80").Value = -4472.0001
$ws.Range("H83").Value = 866.6667
$ws.Range("J83").Value = 866.6667
$ws.Range("L83").Value = 7800.0003
$ws.Range("N83").Value = -17160.0003
$ws.Range("H94").Value = 6432
$ws.Range("I94").Value = 2512
$ws.Range("K94").Value = 7536
$ws.Range("M94").Value = -6860
$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118
$ws.Range("H98").Value = 5963.3335
$ws.Range("I98").Value = 242.5
$ws.Range("J98").Value = 10540
$ws.Range("K98").Value = 727.5
$ws.Range("L98").Value = 31620
$ws.Range("M98").Value = 770.5
$ws.Range("N98").Value = -34616
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("H135").Value = 1576.1111
$ws.Range("J135").Value = 2173.3333
$ws.Range("L135").Value = 19559.9997
$ws.Range("N135").Value = -24629.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488
$ws.Range("H74").Value = 69750
$ws.Range("J74").Value = 69750
$ws.Range("L74").Value = 69750
$ws.Range("N74").Value = -71622
$ws.Range("H77").Value = 69750
$ws.Range("J77").Value = 69750
$ws.Range("L77").Value = 209250
$ws.Range("N77").Value = -218610
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 3977.1914
$ws.Range("I132").Value = 3818.9167
$ws.Range("J132").Value = 4142.3477
$ws.Range("K132").Value = 11456.7501
$ws.Range("L132").Value = 12427.0431
$ws.Range("M132").Value = -8926.750100000001
$ws.Range("N132").Value = -17487.0431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1575
$ws.Range("I61").Value = 1575
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1575
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1373
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 1476.0526
$ws.Range("I82").Value = 834.625
$ws.Range("J82").Value = 1942.5454
$ws.Range("K82").Value = 834.625
$ws.Range("L82").Value = 1942.5454
$ws.Range("M82").Value = -473.625
$ws.Range("N82").Value = -2664.5454
$ws.Range("H85").Value = 1476.0526
$ws.Range("I85").Value = 834.625
$ws.Range("J85").Value = 1942.5454
$ws.Range("K85").Value = 834.625
$ws.Range("L85").Value = 1942.5454
$ws.Range("M85").Value = 413.375
$ws.Range("N85").Value = -4438.5454
$ws.Range("H93").Value = 1148.7646
$ws.Range("I93").Value = 1138.6
$ws.Range("J93").Value = 1163.2858
$ws.Range("K93").Value = 1138.6
$ws.Range("L93").Value = 1163.2858
$ws.Range("M93").Value = 109.4000000000001
$ws.Range("N93").Value = -3659.2858
$ws.Range("H113").Value = 1575
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1575
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 595
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 9234.218999999999
$ws.Range("I132").Value = 3327.889
$ws.Range("K132").Value = 9983.667000000001
$ws.Range("M132").Value = -7453.667000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8933.333000000001
$ws.Range("J74").Value = 8933.333000000001
$ws.Range("L74").Value = 8933.333000000001
$ws.Range("N74").Value = -10805.333
$ws.Range("H77").Value = 8933.333000000001
$ws.Range("J77").Value = 8933.333000000001
$ws.Range("L77").Value = 26799.999
$ws.Range("N77").Value = -36159.999
$ws.Range("H81").Value = 1863.5883
$ws.Range("I81").Value = 898.4167
$ws.Range("J81").Value = 4180
$ws.Range("K81").Value = 1796.8334
$ws.Range("L81").Value = 8360
$ws.Range("M81").Value = -735.8334
$ws.Range("N81").Value = -10482
$ws.Range("H84").Value = 1863.5883
$ws.Range("I84").Value = 898.4167
$ws.Range("J84").Value = 4180
$ws.Range("K84").Value = 8984.166999999999
$ws.Range("L84").Value = 41800
$ws.Range("M84").Value = -3680.166999999999
$ws.Range("N84").Value = -52408
$ws.Range("H132").Value = 1961.3334
$ws.Range("I132").Value = 1253.862
$ws.Range("J132").Value = 3243.625
$ws.Range("K132").Value = 3761.586
$ws.Range("L132").Value = 9730.875
$ws.Range("M132").Value = -1231.586
$ws.Range("N132").Value = -14790.875
